# Quarter rollover: add 2022-Q4 data
#  1. Duplicate the "2022-Q3" sheet to create a new "2022-Q4" sheet placed
#     right before it (i.e. right after "总计"), then overwrite its fund
#     values with the new Q4 numbers.
#  2. Update the "总计" (Total) summary sheet: shift the existing rows down
#     by one and insert the new 2022-Q4 summary row at the top.

$wb = $excel.ActiveWorkbook

# ---- 1. Create the "2022-Q4" worksheet from a copy of "2022-Q3" ----
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# D:G on rows 2-4 hold numeric-looking text (e.g. "0.90") in the source
# workbook, not real numbers. Force text storage (NumberFormat "@") then
# restore the cell style to Normal so no stray number-format style lingers
# in styles.xml once the literal text value has been committed.
$textCells = "D2","E2","F2","G2","D3","E3","F3","G3","D4","E4","F4","G4"
foreach ($addr in $textCells) {
    $q4.Range($addr).NumberFormat = "@"
}

$q4.Range("D2").Value = "1.83"
$q4.Range("E2").Value = "38.98"
$q4.Range("F2").Value = "0.90"
$q4.Range("G2").Value = "0.0165"
$q4.Range("H2").Value = 7

$q4.Range("D3").Value = "1.17"
$q4.Range("E3").Value = "38.98"
$q4.Range("F3").Value = "0.90"
$q4.Range("G3").Value = "0.0105"
$q4.Range("H3").Value = 7

$q4.Range("D4").Value = "0.16"
$q4.Range("E4").Value = "38.98"
$q4.Range("F4").Value = "0.90"
$q4.Range("G4").Value = "0.0014"
$q4.Range("H4").Value = 7

foreach ($addr in $textCells) {
    $q4.Range($addr).Style = "Normal"
}

# ---- 2. Update the "总计" summary sheet ----
$total = $wb.Worksheets.Item("总计")

# Copy the formatting of the last existing data row down into the new row 5
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.01

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.07000000000000001

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q3"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.04

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.03
